$d = $word.ActiveDocument
$find = $d.Content.Find
$find.Text = "Persoana fizica"
$find.Replacement.Text = "Denumire: IDENTITY LEARNING  SRL^pCUI/Tax ID no: 22686237^pAdresa/Adress: JUD. ILFOV, ORŞ. MĂGURELE, STR. CIOCÂRLIEI, NR.11, C2^pRegistrul comertului/Registration no: J23/3344/2017^pEmail: rares.goiceanu@arsek.ro"
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)
